$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 11 ("geo" row) entirely, shifting rows 12/13 up
$ws.Rows(11).Delete()

# Update geo_lon description (now row 11, col C)
$ws.Range("C11").Value = "geographic coordinate: longitude"

# Update geo_lat description (now row 12, col C)
$ws.Range("C12").Value = "geographic coordinate: latitude"
